# EducationalScheduler/sheet.xlsx edit:
# "added todos and got double to long conversion"
#
# Adds a new placeholder/"todo" row (r=3) of "a" values to the
# "Account Information" sheet, mirroring the structure of the existing
# data row, with the last column holding a numeric id
# (89128184 -- the double-to-long conversion fix from the commit msg).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("A3").Value = "a"
$ws1.Range("B3").Value = "a"
$ws1.Range("C3").Value = "a"
$ws1.Range("D3").Value = "a"
$ws1.Range("E3").Value = "a"
$ws1.Range("F3").Value = 89128184
